$wb = $excel.ActiveWorkbook

# --- Rename the "allowed_room" sheet back to its default "Sheet4" name ---
$ws4 = $wb.Worksheets.Item("allowed_room")
$ws4.Name = "Sheet4"

# --- Select A1:D2 on Sheet4 (matches the final on-disk selection state) ---
$ws4.Range("A1:D2").Select()

# --- Add a new worksheet right after Sheet4; it becomes "Sheet5" and the active tab ---
$ws5 = $wb.Worksheets.Add($null, $ws4)

# --- Populate Sheet5 with a copy of Sheet4's header row + first data row ---
$ws5.Range("A1").Value = "block_name"
$ws5.Range("B1").Value = "room_name"
$ws5.Range("C1").Value = "row"
$ws5.Range("D1").Value = "column"

$ws5.Range("A2").Value = "F"
$ws5.Range("B2").Value = "A1"
$ws5.Range("C2").Value = 12
$ws5.Range("D2").Value = 4

# --- Select A1:D2 on the new sheet too, and make sure it is the active sheet/tab ---
$ws5.Range("A1:D2").Select()
$ws5.Activate()
